$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GE")
$ws.Activate()

# Insert two new columns before the existing "expected" column (G), which
# shifts it to column I and shifts the formulas in it automatically.
$ws.Range("G1:H1").EntireColumn.Insert()

# New header labels (order matches how the shared-string table grows in the
# target workbook: flip_angle, echo_time, repetition_time, flip_angle_radians)
$ws.Range("G5").Value = "flip_angle"
$ws.Range("E5").Value = "echo_time"
$ws.Range("F5").Value = "repetition_time"
$ws.Range("H5").Value = "flip_angle_radians"

# Best-effort column widths for the two new columns (engine quantises
# ColumnWidth to 1/6 character-width steps, so this is the closest match).
$ws.Columns.Item(7).ColumnWidth = 9.0
$ws.Columns.Item(8).ColumnWidth = 16.666666666666664

# Flip angle (degrees) for every data row - the new equation supports
# arbitrary flip angles, the validation sheet exercises 90 degrees.
$ws.Range("G6:G16").Value = 90

# Flip angle in radians, used by the new signal equation.
$ws.Range("H6:H16").Formula = "=RADIANS(G6)"

# New Gradient Echo signal equation which allows arbitrary flip angles.
$ws.Range("I6:I16").Formula = "=SIN(H6)*C6*(1-EXP(-F6/A6))/(1-COS(H6)*EXP(-F6/A6)-EXP(-F6/B6)*(EXP(-F6/A6)-COS(H6)))*EXP(-E6/D6)"

# Mirror the author's final selection/active-sheet state.
$ws.Range("I6:I16").Select()
